$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "F2 format:" $ws.Range("F2").NumberFormat
